$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 308, pushing the existing rows 308-314
# (and their data) down to 310-316.
$ws.Rows.Item(308).Resize(2).Insert()

# Row 308 (new): Calera - Repollo, Primera, week of 2021-09-09
$ws.Cells.Item(308, 1).Value = 3
$ws.Cells.Item(308, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(308, 3).Value = "Coquimbo"
$ws.Cells.Item(308, 4).Value = 44448
$ws.Cells.Item(308, 5).Value = 5
$ws.Cells.Item(308, 6).Value = 100112006
$ws.Cells.Item(308, 7).Value = "Repollo"
$ws.Cells.Item(308, 8).Value = "Crespo record"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 3100
$ws.Cells.Item(308, 11).Value = 550
$ws.Cells.Item(308, 12).Value = 600
$ws.Cells.Item(308, 13).Value = 576
$ws.Cells.Item(308, 14).Value = "`$/unidad"
$ws.Cells.Item(308, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(308, 16).Value = 576
$ws.Cells.Item(308, 17).Value = 1
$ws.Cells.Item(308, 18).Value = "Hortaliza"

# Row 309 (new): Calera - Repollo, Segunda, week of 2021-09-09
$ws.Cells.Item(309, 1).Value = 3
$ws.Cells.Item(309, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(309, 3).Value = "Coquimbo"
$ws.Cells.Item(309, 4).Value = 44448
$ws.Cells.Item(309, 5).Value = 5
$ws.Cells.Item(309, 6).Value = 100112006
$ws.Cells.Item(309, 7).Value = "Repollo"
$ws.Cells.Item(309, 8).Value = "Crespo record"
$ws.Cells.Item(309, 9).Value = "Segunda"
$ws.Cells.Item(309, 10).Value = 1900
$ws.Cells.Item(309, 11).Value = 500
$ws.Cells.Item(309, 12).Value = 500
$ws.Cells.Item(309, 13).Value = 500
$ws.Cells.Item(309, 14).Value = "`$/unidad"
$ws.Cells.Item(309, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(309, 16).Value = 500
$ws.Cells.Item(309, 17).Value = 1
$ws.Cells.Item(309, 18).Value = "Hortaliza"
